$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.134.73'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '2.499.41'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.539'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.86'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.92'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0804'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').Value = '2.890.52'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '2.515.55'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.833'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '47.970.43'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').Value = '0.0₃0934'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '271.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.73'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.96%  '
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0778'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.94'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.59'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.88'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '120.97'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.06'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0304'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').Value = '2.007.45'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.12%  '
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.70'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.78%  '
